$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column F (Posted At) to text so date-like strings are not
# auto-converted to date serials, then clear the format afterwards so
# no stray style index is left on the cells.
$ws.Range("F2:F26").NumberFormat = "@"

$ws.Cells.Item(2, 1).Value = 'Senior Data and Software Engineer I (API) (32_2026.1)'
$ws.Cells.Item(2, 2).Value = 'Affinity Solutions'
$ws.Cells.Item(2, 3).Value = 'New York, NY, US USA'
$ws.Cells.Item(2, 4).Value = 26.7
$ws.Cells.Item(2, 5).Value = 'Generative AI, RAG, Prompt Engineering, Cortex, TensorFlow, S3, Glue, Athena, Redshift, Data Lake'
$ws.Cells.Item(2, 6).Value = '2026-02-24'
$ws.Cells.Item(2, 7).Value = 'https://www.indeed.com/viewjob?jk=6ced5764d3537ea0'

$ws.Cells.Item(3, 1).Value = 'Sr. Data Quality Engineer I (37_2026.1)'
$ws.Cells.Item(3, 2).Value = 'Affinity Solutions'
$ws.Cells.Item(3, 3).Value = 'New York, NY, US USA'
$ws.Cells.Item(3, 4).Value = 25.6
$ws.Cells.Item(3, 5).Value = 'RAG, S3, Glue, Athena, Redshift, BigQuery, Data Lake, Apache Airflow, CI/CD, Jenkins'
$ws.Cells.Item(3, 6).Value = '2026-02-24'
$ws.Cells.Item(3, 7).Value = 'https://www.indeed.com/viewjob?jk=c7933ce6e87f7544'

$ws.Cells.Item(4, 1).Value = 'Identity AI / ML Engineer'
$ws.Cells.Item(4, 2).Value = 'Openkyber'
$ws.Cells.Item(4, 3).Value = 'AK, US USA'
$ws.Cells.Item(4, 4).Value = 17.8
$ws.Cells.Item(4, 5).Value = 'AI Engineer, Generative AI, LangChain, RAG, Hugging Face, FAISS, Pinecone, Prompt Engineering, TensorFlow, PyTorch'
$ws.Cells.Item(4, 6).Value = '2026-02-24'
$ws.Cells.Item(4, 7).Value = 'https://www.indeed.com/viewjob?jk=91786fa9f4c64147'

$ws.Cells.Item(5, 1).Value = 'AI Technical Architect'
$ws.Cells.Item(5, 2).Value = 'Capco'
$ws.Cells.Item(5, 3).Value = 'Orlando, FL, US USA'
$ws.Cells.Item(5, 4).Value = 17.8
$ws.Cells.Item(5, 5).Value = 'AI Engineer, Data Scientist, Generative AI, LangChain, RAG, FAISS, Pinecone, Prompt Engineering, TensorFlow, PyTorch'
$ws.Cells.Item(5, 6).Value = '2026-02-24'
$ws.Cells.Item(5, 7).Value = 'https://www.indeed.com/viewjob?jk=42acc1d7c9836f52'

$ws.Cells.Item(6, 1).Value = 'AI Technical Architect'
$ws.Cells.Item(6, 2).Value = 'Capco'
$ws.Cells.Item(6, 3).Value = 'New York, NY, US USA'
$ws.Cells.Item(6, 4).Value = 17.8
$ws.Cells.Item(6, 5).Value = 'AI Engineer, Data Scientist, Generative AI, LangChain, RAG, FAISS, Pinecone, Prompt Engineering, TensorFlow, PyTorch'
$ws.Cells.Item(6, 6).Value = '2026-02-24'
$ws.Cells.Item(6, 7).Value = 'https://www.indeed.com/viewjob?jk=06e9d21d7dc42caa'

$ws.Cells.Item(7, 1).Value = 'Software Engineer III - Full Stack + AWS + Elastic / Open Search'
$ws.Cells.Item(7, 2).Value = 'JPMorganChase'
$ws.Cells.Item(7, 3).Value = 'Plano, TX, US USA'
$ws.Cells.Item(7, 4).Value = 16.7
$ws.Cells.Item(7, 5).Value = 'RAG, S3, EC2, CI/CD, Jenkins, GitHub Actions, Terraform, Git, Kafka, PostgreSQL'
$ws.Cells.Item(7, 6).Value = '2026-02-24'
$ws.Cells.Item(7, 7).Value = 'https://www.indeed.com/viewjob?jk=9d9fbf9321e91c82'

$ws.Cells.Item(8, 1).Value = 'Data Architect Senior'
$ws.Cells.Item(8, 2).Value = 'The University of Michigan'
$ws.Cells.Item(8, 3).Value = 'Ann Arbor, MI, US USA'
$ws.Cells.Item(8, 4).Value = 15.6
$ws.Cells.Item(8, 5).Value = 'RAG, Synapse, CI/CD, Git, Snowflake, Databricks, PostgreSQL, Tableau, Power BI, Python'
$ws.Cells.Item(8, 6).Value = '2026-02-24'
$ws.Cells.Item(8, 7).Value = 'https://www.indeed.com/viewjob?jk=4ccb9adf6566f3d5'

$ws.Cells.Item(9, 1).Value = 'Data Engineer III'
$ws.Cells.Item(9, 2).Value = 'JPMorganChase'
$ws.Cells.Item(9, 3).Value = 'Columbus, OH, US USA'
$ws.Cells.Item(9, 4).Value = 15.6
$ws.Cells.Item(9, 5).Value = 'RAG, Glue, Athena, Kubernetes, CI/CD, Terraform, Git, Snowflake, PySpark, Kafka'
$ws.Cells.Item(9, 6).Value = '2026-02-24'
$ws.Cells.Item(9, 7).Value = 'https://www.indeed.com/viewjob?jk=b3e5e0e8fb1f2fdf'

$ws.Cells.Item(10, 1).Value = 'Identity AI / ML Engineer'
$ws.Cells.Item(10, 2).Value = 'Openkyber'
$ws.Cells.Item(10, 3).Value = 'AK, US USA'
$ws.Cells.Item(10, 4).Value = 14.4
$ws.Cells.Item(10, 5).Value = 'Data Scientist, Machine Learning Engineer, TensorFlow, PyTorch, AWS SageMaker, Azure ML, Docker, Kubernetes, CI/CD, Python'
$ws.Cells.Item(10, 6).Value = '2026-02-24'
$ws.Cells.Item(10, 7).Value = 'https://www.indeed.com/viewjob?jk=21d02cf905324c4e'

$ws.Cells.Item(11, 1).Value = 'Data Scientist - Remote'
$ws.Cells.Item(11, 2).Value = 'Ocean Blue Solutions'
$ws.Cells.Item(11, 3).Value = 'Columbus, OH, US USA'
$ws.Cells.Item(11, 4).Value = 14.4
$ws.Cells.Item(11, 5).Value = 'Data Scientist, RAG, TensorFlow, PyTorch, PySpark, Tableau, Power BI, Matplotlib, Python, R'
$ws.Cells.Item(11, 6).Value = '2026-02-24'
$ws.Cells.Item(11, 7).Value = 'https://www.indeed.com/viewjob?jk=244c2e160f4e3ba6'

$ws.Cells.Item(12, 1).Value = 'Identity Shield & Consumer Fraud - FDP (Backend)'
$ws.Cells.Item(12, 2).Value = 'Ally Financial'
$ws.Cells.Item(12, 3).Value = 'Charlotte, NC, US USA'
$ws.Cells.Item(12, 4).Value = 14.4
$ws.Cells.Item(12, 5).Value = 'Data Scientist, RAG, Glue, Kinesis, Terraform, Git, PostgreSQL, NoSQL, Python, SQL'
$ws.Cells.Item(12, 6).Value = '2025-11-14'
$ws.Cells.Item(12, 7).Value = 'https://www.indeed.com/viewjob?jk=a94b094102920ddc'

$ws.Cells.Item(13, 1).Value = 'Sr ML Ops Engineer'
$ws.Cells.Item(13, 2).Value = 'Early Warning Services'
$ws.Cells.Item(13, 3).Value = 'San Francisco, CA, US USA'
$ws.Cells.Item(13, 4).Value = 13.3
$ws.Cells.Item(13, 5).Value = 'Data Scientist, RAG, MLflow, Docker, Kubernetes, CI/CD, Git, Hadoop, Python, R'
$ws.Cells.Item(13, 6).Value = '2026-02-24'
$ws.Cells.Item(13, 7).Value = 'https://www.indeed.com/viewjob?jk=daf4a4241dc4709a'

$ws.Cells.Item(14, 1).Value = 'Software Engineer'
$ws.Cells.Item(14, 2).Value = 'Anaplan'
$ws.Cells.Item(14, 3).Value = 'Philadelphia, PA, US USA'
$ws.Cells.Item(14, 4).Value = 12.2
$ws.Cells.Item(14, 5).Value = 'Docker, Kubernetes, CI/CD, Git, PostgreSQL, NoSQL, SQL, R, Java, Scala'
$ws.Cells.Item(14, 6).Value = '2026-02-24'
$ws.Cells.Item(14, 7).Value = 'https://www.indeed.com/viewjob?jk=d985dfab97c4672d'

$ws.Cells.Item(15, 1).Value = '(1357) Senior Reliability Engineer'
$ws.Cells.Item(15, 2).Value = 'Nearsure'
$ws.Cells.Item(15, 3).Value = 'Remote, US USA'
$ws.Cells.Item(15, 4).Value = 12.2
$ws.Cells.Item(15, 5).Value = 'Generative AI, EC2, Kubernetes, AKS, CI/CD, GitHub Actions, Terraform, Git, R, Scala'
$ws.Cells.Item(15, 6).Value = '2026-02-24'
$ws.Cells.Item(15, 7).Value = 'https://www.indeed.com/viewjob?jk=b3c8f5bbb025388b'

$ws.Cells.Item(16, 1).Value = 'Machine Learning Engineer III - Generative AI - Windreich Department of Artificial Intelligence & Human Health Research'
$ws.Cells.Item(16, 2).Value = 'Mount Sinai Health System'
$ws.Cells.Item(16, 3).Value = 'New York, NY, US USA'
$ws.Cells.Item(16, 4).Value = 11.1
$ws.Cells.Item(16, 5).Value = 'Data Scientist, Machine Learning Engineer, Generative AI, RAG, Jenkins, Terraform, Git, NoSQL, SQL, R'
$ws.Cells.Item(16, 6).Value = '2026-02-24'
$ws.Cells.Item(16, 7).Value = 'https://www.indeed.com/viewjob?jk=6c988cf5e8e74111'

$ws.Cells.Item(17, 1).Value = 'Identity AI / ML Engineer'
$ws.Cells.Item(17, 2).Value = 'Openkyber'
$ws.Cells.Item(17, 3).Value = 'AK, US USA'
$ws.Cells.Item(17, 4).Value = 11.1
$ws.Cells.Item(17, 5).Value = 'Data Scientist, Machine Learning Engineer, TensorFlow, PyTorch, AWS SageMaker, Azure ML, Python, R, Scala, A/B Testing'
$ws.Cells.Item(17, 6).Value = '2026-02-24'
$ws.Cells.Item(17, 7).Value = 'https://www.indeed.com/viewjob?jk=144e6727ca38d104'

$ws.Cells.Item(18, 1).Value = 'Identity AI / ML Engineer'
$ws.Cells.Item(18, 2).Value = 'Openkyber'
$ws.Cells.Item(18, 3).Value = 'AK, US USA'
$ws.Cells.Item(18, 4).Value = 11.1
$ws.Cells.Item(18, 5).Value = 'RAG, spaCy, NLTK, S3, Docker, CI/CD, Git, Python, SQL, R'
$ws.Cells.Item(18, 6).Value = '2026-02-24'
$ws.Cells.Item(18, 7).Value = 'https://www.indeed.com/viewjob?jk=65ccf46dbfdb26c1'

$ws.Cells.Item(19, 1).Value = 'AIML Services- Data Platform Engineer'
$ws.Cells.Item(19, 2).Value = 'GSK'
$ws.Cells.Item(19, 3).Value = 'San Francisco, CA, US USA'
$ws.Cells.Item(19, 4).Value = 11.1
$ws.Cells.Item(19, 5).Value = 'Machine Learning Engineer, MLflow, Docker, Kubernetes, CI/CD, Jenkins, Git, Python, R, Scala'
$ws.Cells.Item(19, 6).Value = '2026-02-24'
$ws.Cells.Item(19, 7).Value = 'https://www.indeed.com/viewjob?jk=a25ec63d2a45d72c'

$ws.Cells.Item(20, 1).Value = 'AIML Services- Data Platform Engineer'
$ws.Cells.Item(20, 2).Value = 'GSK'
$ws.Cells.Item(20, 3).Value = 'Seattle, WA, US USA'
$ws.Cells.Item(20, 4).Value = 11.1
$ws.Cells.Item(20, 5).Value = 'Machine Learning Engineer, MLflow, Docker, Kubernetes, CI/CD, Jenkins, Git, Python, R, Scala'
$ws.Cells.Item(20, 6).Value = '2026-02-24'
$ws.Cells.Item(20, 7).Value = 'https://www.indeed.com/viewjob?jk=41cbea8637a10d68'

$ws.Cells.Item(21, 1).Value = 'AIML Services- Data Platform Engineer'
$ws.Cells.Item(21, 2).Value = 'GSK'
$ws.Cells.Item(21, 3).Value = 'Cambridge, MA, US USA'
$ws.Cells.Item(21, 4).Value = 11.1
$ws.Cells.Item(21, 5).Value = 'Machine Learning Engineer, MLflow, Docker, Kubernetes, CI/CD, Jenkins, Git, Python, R, Scala'
$ws.Cells.Item(21, 6).Value = '2026-02-24'
$ws.Cells.Item(21, 7).Value = 'https://www.indeed.com/viewjob?jk=94ad64f2466c8d7d'

$ws.Cells.Item(22, 1).Value = 'Senior Data Analyst'
$ws.Cells.Item(22, 2).Value = 'HDR'
$ws.Cells.Item(22, 3).Value = 'Charleston, WV, US USA'
$ws.Cells.Item(22, 4).Value = 11.1
$ws.Cells.Item(22, 5).Value = 'Data Scientist, RAG, Snowflake, Tableau, Power BI, Matplotlib, Seaborn, Python, SQL, R'
$ws.Cells.Item(22, 6).Value = '2026-02-24'
$ws.Cells.Item(22, 7).Value = 'https://www.indeed.com/viewjob?jk=30f4d6d3f413c72c'

$ws.Cells.Item(23, 1).Value = 'Software Engineer II - DevOps, Platform Engineering'
$ws.Cells.Item(23, 2).Value = 'JPMorganChase'
$ws.Cells.Item(23, 3).Value = 'Chicago, IL, US USA'
$ws.Cells.Item(23, 4).Value = 10
$ws.Cells.Item(23, 5).Value = 'RAG, Kubernetes, CI/CD, Jenkins, Terraform, Git, Python, R, Scala'
$ws.Cells.Item(23, 6).Value = '2026-02-24'
$ws.Cells.Item(23, 7).Value = 'https://www.indeed.com/viewjob?jk=11202771fd209720'

$ws.Cells.Item(24, 1).Value = 'AI Engineer'
$ws.Cells.Item(24, 2).Value = 'Capco'
$ws.Cells.Item(24, 3).Value = 'New York, NY, US USA'
$ws.Cells.Item(24, 4).Value = 10
$ws.Cells.Item(24, 5).Value = 'AI Engineer, Generative AI, RAG, CI/CD, Git, Python, R, Java, Scala'
$ws.Cells.Item(24, 6).Value = '2026-02-24'
$ws.Cells.Item(24, 7).Value = 'https://www.indeed.com/viewjob?jk=5b9f8300c935d7e1'

$ws.Cells.Item(25, 1).Value = 'AI Engineer'
$ws.Cells.Item(25, 2).Value = 'Capco'
$ws.Cells.Item(25, 3).Value = 'Orlando, FL, US USA'
$ws.Cells.Item(25, 4).Value = 10
$ws.Cells.Item(25, 5).Value = 'AI Engineer, Generative AI, RAG, CI/CD, Git, Python, R, Java, Scala'
$ws.Cells.Item(25, 6).Value = '2026-02-24'
$ws.Cells.Item(25, 7).Value = 'https://www.indeed.com/viewjob?jk=22b8ff99671f9ffa'

$ws.Cells.Item(26, 1).Value = 'AI Engineer'
$ws.Cells.Item(26, 2).Value = 'Capco'
$ws.Cells.Item(26, 3).Value = 'Dallas, TX, US USA'
$ws.Cells.Item(26, 4).Value = 10
$ws.Cells.Item(26, 5).Value = 'AI Engineer, Generative AI, RAG, CI/CD, Git, Python, R, Java, Scala'
$ws.Cells.Item(26, 6).Value = '2026-02-24'
$ws.Cells.Item(26, 7).Value = 'https://www.indeed.com/viewjob?jk=1f18056c1a09ca54'

$ws.Range("F2:F26").ClearFormats()
